$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("harp expander v1")

# Row 5 is the R1 resistor entry. Update its MPN, distributor part number,
# and value to reflect the new 910 ohm part (was 330 ohm).
$ws.Range("E5").Value = "RK73B1ETTP911J"
$ws.Range("H5").Value = "2019-RK73B1ETTP911JDKR-ND"
$ws.Range("B5").Value = "910 1/10W 5%"

# Reflect the saved cursor/selection position from the author's session.
$ws.Range("D19").Select()
